$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 283.53333
$ws.Range("I2").Value = 273.5
$ws.Range("J2").Value = 424
$ws.Range("K2").Value = 273.5
$ws.Range("L2").Value = 424
$ws.Range("M2").Value = -160.5
$ws.Range("N2").Value = -650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3151.25
$ws.Range("I2").Value = 837.2
$ws.Range("J2").Value = 7008
$ws.Range("K2").Value = 837.2
$ws.Range("L2").Value = 7008
$ws.Range("M2").Value = -724.2
$ws.Range("N2").Value = -7234

$ws.Range("H36").Value = 35000
$ws.Range("I36").Value = 35000
$ws.Range("K36").Value = 35000
$ws.Range("M36").Value = -34654

$ws.Range("H61").Value = 5498.25
$ws.Range("I61").Value = 5497.5
$ws.Range("K61").Value = 5497.5
$ws.Range("M61").Value = -5285.5

$ws.Range("H92").Value = 50550
$ws.Range("J92").Value = 50550
$ws.Range("L92").Value = 50550
$ws.Range("N92").Value = -55542

$ws.Range("H116").Value = 3151.25
$ws.Range("I116").Value = 837.2
$ws.Range("J116").Value = 7008
$ws.Range("K116").Value = 837.2
$ws.Range("L116").Value = 7008
$ws.Range("M116").Value = 1456.8
$ws.Range("N116").Value = -11596

$ws.Range("H122").Value = 1498.2142
$ws.Range("I122").Value = 1411.1
$ws.Range("J122").Value = 1716
$ws.Range("K122").Value = 4233.299999999999
$ws.Range("L122").Value = 5148
$ws.Range("M122").Value = -1783.299999999999
$ws.Range("N122").Value = -10048

$ws.Range("H136").Value = 5498.25
$ws.Range("I136").Value = 5497.5
$ws.Range("K136").Value = 16492.5
$ws.Range("M136").Value = -13942.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3151.25
$ws.Range("I3").Value = 837.2
$ws.Range("J3").Value = 7008
$ws.Range("K3").Value = 837.2
$ws.Range("L3").Value = 7008
$ws.Range("M3").Value = -723.2
$ws.Range("N3").Value = -7236

$ws.Range("H7").Value = 6333950
$ws.Range("I7").Value = 6333950
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 6333950
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -6333837
$ws.Range("N7").ClearContents()

$ws.Range("H10").Value = 1493
$ws.Range("J10").Value = 1493
$ws.Range("L10").Value = 1493
$ws.Range("N10").Value = -1773

$ws.Range("H107").Value = 668.55554
$ws.Range("I107").Value = 578.6667
$ws.Range("J107").Value = 848.3333
$ws.Range("K107").Value = 578.6667
$ws.Range("L107").Value = 848.3333
$ws.Range("M107").Value = 1341.3333
$ws.Range("N107").Value = -4688.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H62").Value = 1500
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 1500
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -16240

$ws.Range("H141").Value = 398355.16
$ws.Range("J141").Value = 398355.16
$ws.Range("L141").Value = 398355.16
$ws.Range("N141").Value = -408715.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 498.5
$ws.Range("I3").Value = 498.5
$ws.Range("K3").Value = 1495.5
$ws.Range("M3").Value = -1383.5

$ws.Range("H41").Value = 856
$ws.Range("I41").Value = 335.42856
$ws.Range("J41").Value = 4500
$ws.Range("K41").Value = 1006.28568
$ws.Range("L41").Value = 13500
$ws.Range("M41").Value = -668.28568
$ws.Range("N41").Value = -14176

$ws.Range("H122").Value = 1475
$ws.Range("J122").Value = 1483
$ws.Range("L122").Value = 13347
$ws.Range("N122").Value = -18247

$ws.Range("H131").Value = 894.25
$ws.Range("J131").Value = 930.24243
$ws.Range("L131").Value = 2790.72729
$ws.Range("N131").Value = -12870.72729

$ws.Range("H138").Value = 3299.8
$ws.Range("I138").Value = 3832.6667
$ws.Range("K138").Value = 11498.0001
$ws.Range("M138").Value = -6358.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 746.875
$ws.Range("I107").Value = 567.8570999999999
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 567.8570999999999
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1352.1429
$ws.Range("N107").Value = -5840

$ws.Range("H122").Value = 1999.6
$ws.Range("I122").Value = 1857
$ws.Range("K122").Value = 5571
$ws.Range("M122").Value = -3121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 905.8461
$ws.Range("I22").Value = 909
$ws.Range("J22").Value = 898.75
$ws.Range("K22").Value = 909
$ws.Range("L22").Value = 898.75
$ws.Range("M22").Value = -614
$ws.Range("N22").Value = -1488.75

$ws.Range("H27").Value = 905.8461
$ws.Range("I27").Value = 909
$ws.Range("J27").Value = 898.75
$ws.Range("K27").Value = 909
$ws.Range("L27").Value = 898.75
$ws.Range("M27").Value = -802
$ws.Range("N27").Value = -1112.75

$ws.Range("H140").Value = 79000
$ws.Range("J140").Value = 79000
$ws.Range("L140").Value = 79000
$ws.Range("N140").Value = -89360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 36913.332
$ws.Range("I112").Value = 10000
$ws.Range("J112").Value = 50370
$ws.Range("K112").Value = 10000
$ws.Range("L112").Value = 50370
$ws.Range("M112").Value = -8523
$ws.Range("N112").Value = -53324

$ws.Range("H122").Value = 486.9091
$ws.Range("I122").Value = 493.1
$ws.Range("J122").Value = 425
$ws.Range("K122").Value = 1479.3
$ws.Range("L122").Value = 1275
$ws.Range("M122").Value = 970.6999999999998
$ws.Range("N122").Value = -6175

$ws.Range("H125").Value = 62500
$ws.Range("J125").Value = 62500
$ws.Range("L125").Value = 62500
$ws.Range("N125").Value = -72340
